$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("n1_d40")
$ws.Range("F2").Value = 0.0615436
$ws.Range("G2").Value = 8.699999999999999
$ws.Range("F3").Value = 0.0625984
$ws.Range("G3").Value = 8.800000000000001
$ws.Range("F4").Value = 0.061477
$ws.Range("G4").Value = 8.699999999999999
$ws.Range("F5").Value = 0.061504
$ws.Range("G5").Value = 8.800000000000001
$ws.Range("F6").Value = 0.0615597
$ws.Range("G6").Value = 8.6
$ws.Range("F7").Value = 0.0614529
$ws.Range("G7").Value = 8.699999999999999
$ws.Range("F8").Value = 0.0615633
$ws.Range("G8").Value = 8.699999999999999
$ws.Range("F9").Value = 0.0626086
$ws.Range("G9").Value = 9.5
$ws.Range("F10").Value = 0.0625894
$ws.Range("G10").Value = 9.5
$ws.Range("F11").Value = 0.0626419
$ws.Range("G11").Value = 9.5
$ws.Range("F12").Value = 0.06195388
$ws.Range("G12").Value = 8.949999999999999

$ws = $wb.Worksheets.Item("n1_d60")
$ws.Range("F2").Value = 0.108558
$ws.Range("G2").Value = 15.7
$ws.Range("F3").Value = 0.107557
$ws.Range("G3").Value = 15.6
$ws.Range("F4").Value = 0.107399
$ws.Range("G4").Value = 15.7
$ws.Range("F5").Value = 0.110661
$ws.Range("G5").Value = 15.7
$ws.Range("F6").Value = 0.108539
$ws.Range("G6").Value = 15.7
$ws.Range("F7").Value = 0.107461
$ws.Range("G7").Value = 15.7
$ws.Range("F8").Value = 0.109457
$ws.Range("G8").Value = 15.7
$ws.Range("F9").Value = 0.109663
$ws.Range("G9").Value = 15.7
$ws.Range("F10").Value = 0.109648
$ws.Range("G10").Value = 15.7
$ws.Range("F11").Value = 0.107345
$ws.Range("G11").Value = 15.7
$ws.Range("F12").Value = 0.1086288
$ws.Range("G12").Value = 15.69

$ws = $wb.Worksheets.Item("n1_d80")
$ws.Range("F2").Value = 0.155685
$ws.Range("G2").Value = 23
$ws.Range("F3").Value = 0.156682
$ws.Range("G3").Value = 23.6
$ws.Range("F4").Value = 0.156595
$ws.Range("G4").Value = 23
$ws.Range("F5").Value = 0.154424
$ws.Range("G5").Value = 23.6
$ws.Range("F6").Value = 0.15968
$ws.Range("G6").Value = 23.2
$ws.Range("F7").Value = 0.155954
$ws.Range("G7").Value = 23.4
$ws.Range("F8").Value = 0.156764
$ws.Range("G8").Value = 23
$ws.Range("F9").Value = 0.157574
$ws.Range("G9").Value = 23.3
$ws.Range("F10").Value = 0.156624
$ws.Range("G10").Value = 23
$ws.Range("F11").Value = 0.160933
$ws.Range("G11").Value = 23.2
$ws.Range("F12").Value = 0.1570915
$ws.Range("G12").Value = 23.23

$ws = $wb.Worksheets.Item("n1_d100")
$ws.Range("F2").Value = 0.203614
$ws.Range("G2").Value = 27.1
$ws.Range("F3").Value = 0.201235
$ws.Range("G3").Value = 27.7
$ws.Range("F4").Value = 0.20126
$ws.Range("G4").Value = 27.6
$ws.Range("F5").Value = 0.200355
$ws.Range("G5").Value = 28.5
$ws.Range("F6").Value = 0.203542
$ws.Range("G6").Value = 28.9
$ws.Range("F7").Value = 0.201133
$ws.Range("G7").Value = 27.8
$ws.Range("F8").Value = 0.202215
$ws.Range("G8").Value = 27.5
$ws.Range("F9").Value = 0.204395
$ws.Range("G9").Value = 27.9
$ws.Range("F10").Value = 0.201762
$ws.Range("G10").Value = 28.3
$ws.Range("F11").Value = 0.203813
$ws.Range("G11").Value = 27.9
$ws.Range("F12").Value = 0.2023324
$ws.Range("G12").Value = 27.92

